$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Saturday hours for the week of row 8 (G8) from 6.5 to 8.5.
# This will automatically recalculate the dependent formulas in I8 (row total)
# and I19 (grand total).
$ws.Range("G8").Value = 8.5

# Recalculate to ensure formula results (I8, I19) are refreshed.
$excel.Calculate()

# Update the active cell selection as it appears in the saved file.
$ws.Range("L6").Select()
